$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 2 with new test values (user name / password)
$ws.Range("B2").Value = "Test User name"
$ws.Range("C2").Value = "Test Password"

# Add a new Row 3 with second set of test data
$ws.Range("A3").Value = "TestDB2"
$ws.Range("B3").Value = "Test User name2"
$ws.Range("C3").Value = "Test Password2"
$ws.Range("D3").Value = "localhost2"
$ws.Range("E3").Value = 15212

# Update the active selection to match the new editing location
$ws.Range("G9").Select()
